$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 / 17 swap: PropGW (was row16) <-> DOC_gw (was row17)
$ws.Range("A16").Value = "DOC_gw"
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = "g/m3"

$ws.Range("A17").Value = "PropGW"
$ws.Range("B17").Value = 0.19
$ws.Range("C17").Value = "unitless"

# Row 18 (DOC_sw) contents removed entirely
$ws.Range("A18:C18").ClearContents()

# Update the active selection to C17
$ws.Range("C17").Select()
